$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Split the "The IGB and Plant Biology staff ... mine." / "The Plant
#    Biology Association ..." paragraph into two separate paragraphs,
#    removing the lone-space run that used to join them.
# ------------------------------------------------------------------
$d.Content.Find.Execute( `
    "mine. The Plant Biology Association", $true, $false, $false, $false, $false, `
    $true, 1, $false, "mine.^pThe Plant Biology Association", 2) | Out-Null

# ------------------------------------------------------------------
# 2) Split the "My advisor ... wrong." / "Coauthors ..." paragraph into
#    two separate paragraphs, removing the lone-space run that used to
#    join them.
# ------------------------------------------------------------------
$d.Content.Find.Execute( `
    "wrong. Coauthors", $true, $false, $false, $false, $false, `
    $true, 1, $false, "wrong.^pCoauthors", 2) | Out-Null

# ------------------------------------------------------------------
# 3) Insert a brand-new paragraph ahead of "My advisor ..." that reads
#    "I am especially grateful for consistent, dedicated help from:"
#    with "especially" in italics.
# ------------------------------------------------------------------
$locate = $d.Content
$locate.Find.Execute("My advisor, Evan DeLucia", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$locate.Collapse(1)
$advisorParaIndex = $locate.Paragraphs.Item(1).Index

$advisorP = $d.Paragraphs.Item($advisorParaIndex)
$advisorP.Range.InsertParagraphBefore()

$newP = $d.Paragraphs.Item($advisorParaIndex)
$newStart = $newP.Range.Start

$newRange = $d.Range($newStart, $newStart)
$newRange.InsertAfter("I am especially grateful for consistent, dedicated help from:")

# italicize just the word "especially" inside the new paragraph
$newP2 = $d.Paragraphs.Item($advisorParaIndex)
$italicRng = $newP2.Range
$italicRng.Find.Execute("especially", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$italicRng.Font.Italic = 1

# ------------------------------------------------------------------
# 4) Trim the trailing "And I am especially grateful for consistent,
#    dedicated help from:" clause off the first body paragraph -- that
#    sentiment now lives in its own paragraph (see step 3).
# ------------------------------------------------------------------
$d.Content.Find.Execute( `
    " And I am especially grateful for consistent, dedicated help from:", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null
